$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Add the new "Abstract Title" paragraph style (based on Normal,
#    followed by the existing "Abstract" style).
# ---------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = "345A8A"

# ---------------------------------------------------------------------
# 2) Existing "Abstract" style: tighten the space-before from 300 to
#    100 twips (space-after stays 300).
# ---------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# ---------------------------------------------------------------------
# 3) Add the new "Footnote Block Text" paragraph style (based on
#    "Footnote Text", followed by "Footnote Text").
# ---------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "styles updated"
